$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.193.33"
$ws.Range("E2").Value = "  -1.90%  "

$ws.Range("D3").Value = "1.871.69"
$ws.Range("E3").Value = "  -1.54%  "

$ws.Range("D4").Value = "'0.9999"

$ws.Range("D5").Value = "'306.55"
$ws.Range("E5").Value = "  -1.75%  "

$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").Value = "'0.5198"
$ws.Range("E7").Value = "  -0.59%  "

$ws.Range("D8").Value = "'0.3743"
$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("D9").Value = "'0.07163"
$ws.Range("E9").Value = "  -1.17%  "

$ws.Range("D10").Value = "'0.8934"
$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("D11").Value = "'20.78"
$ws.Range("E11").Value = "  -2.11%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07536"
$ws.Range("E12").Value = "  -1.43%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.869.69"
$ws.Range("E13").Value = "  -1.75%  "

$ws.Range("D14").Value = "'5.310"
$ws.Range("E14").Value = "  -2.51%  "

$ws.Range("D15").Value = "'90.57"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").Value = "'1.0000"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").Value = "'0.000008506"
$ws.Range("E17").Value = "  -2.21%  "

$ws.Range("E18").Value = "  -2.42%  "

$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").Value = "27.220.65"
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("D21").Value = "'5.009"
$ws.Range("E21").Value = "  -2.57%  "

$ws.Range("D22").Value = "2.108.19"
$ws.Range("E22").Value = "  -1.66%  "

$ws.Range("D23").Value = "'10.47"
$ws.Range("E23").Value = "  -3.55%  "

$ws.Range("D24").Value = "'6.478"
$ws.Range("E24").Value = "  -2.10%  "

$ws.Range("D25").Value = "'1.834"
$ws.Range("E25").Value = "  -1.79%  "

$ws.Range("D26").Value = "'146.10"
$ws.Range("E26").Value = "  -4.60%  "

$ws.Range("D27").Value = "'18.00"
$ws.Range("E27").Value = "  -1.69%  "

$ws.Range("D28").Value = "'2.091"
$ws.Range("E28").Value = "  -3.24%  "

$ws.Range("D29").Value = "'113.31"
$ws.Range("E29").Value = "  -1.07%  "

$ws.Range("D30").Value = "'4.666"
$ws.Range("E30").Value = "  -3.69%  "

$ws.Range("D31").Value = "'4.689"
$ws.Range("E31").Value = "  -3.04%  "

$ws.Range("D32").Value = "'0.09279"
$ws.Range("E32").Value = "  +2.13%  "

$ws.Range("D33").Value = "'0.05143"
$ws.Range("E33").Value = "  -2.75%  "

$ws.Range("D34").Value = "'3.083"
$ws.Range("E34").Value = "  -3.27%  "

$ws.Range("D35").Value = "'1.162"
$ws.Range("E35").Value = "  -5.12%  "

$ws.Range("D36").Value = "'0.7285"
$ws.Range("E36").Value = "  -6.52%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02034"
$ws.Range("E37").Value = "  -2.78%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'3.120"
$ws.Range("E38").Value = "  +1.75%  "

$ws.Range("D39").Value = "'2.519"
$ws.Range("E39").Value = "  -2.23%  "

$ws.Range("E40").Value = "  -1.70%  "

$ws.Range("D41").Value = "'0.5324"
$ws.Range("E41").Value = "  -4.04%  "

$ws.Range("D42").Value = "'6.530"
$ws.Range("E42").Value = "  -2.96%  "

$ws.Range("D43").Value = "'116.45"
$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("D44").Value = "'8.365"
$ws.Range("E44").Value = "  -1.73%  "

$ws.Range("D45").Value = "'0.1477"
$ws.Range("E45").Value = "  -2.66%  "

$ws.Range("D46").Value = "'0.4636"
$ws.Range("E46").Value = "  -3.70%  "

$ws.Range("D47").Value = "'0.9997"
$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("E48").Value = "  -4.47%  "

$ws.Range("D49").Value = "'1.566"
$ws.Range("E49").Value = "  -2.93%  "

$ws.Range("D50").Value = "'36.76"
$ws.Range("E50").Value = "  -0.70%  "

$ws.Range("E51").Value = "  -4.62%  "
